# Generate Report for Handback
# Update "Correspond Handoff Datetime" (col E) and "Correspond Handback DateTime" (col H)
# for the second data row (row 3) on both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 06:33:46"
$wsZhCn.Range("H3").Value = "2016-03-22 06:34:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 06:33:56"
$wsDeDe.Range("H3").Value = "2016-03-22 06:34:40"
